$wb = $excel.ActiveWorkbook

# --- production_rates: remove negative sign from production rates, select column C ---
$wsProd = $wb.Worksheets.Item("production_rates")
$wsProd.Activate() | Out-Null
$prodRates = $wsProd.Range("B2:B17")
foreach ($cell in $prodRates.Cells) {
    $cell.Value = -1 * $cell.Value()
}
$wsProd.Range("C:C").Select() | Out-Null

# --- degradation_rates: remove negative sign from degradation rates, select column C ---
$wsDeg = $wb.Worksheets.Item("degradation_rates")
$wsDeg.Activate() | Out-Null
$degRates = $wsDeg.Range("B2:B17")
foreach ($cell in $degRates.Cells) {
    $cell.Value = -1 * $cell.Value()
}
$wsDeg.Range("C:C").Select() | Out-Null
